$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 73 (2021-11-07): After Midnight Part 1 further exploration
$ws.Range("B73").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C73").Value = 0.15
$ws.Range("D73").Value = "only anki"

# Row 74 (2021-11-08): further exploration
$ws.Range("B74").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C74").Value = 0.75
$ws.Range("D74").Value = "further exploration and 1 small problem"

$ws.Range("D80").Formula = "=SUM(C73:C79)"
$ws.Range("C81").Formula = "=SUBTOTAL(109,Table1[Hours])"

# Mirror the final view/selection state from the edit session
$ws.Range("C67").Select()
$ws.Range("D75").Select()

$wb.Save()
